$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The date lives in the second paragraph of the subtitle placeholder:
# "Novemeber 20, 2024" -> split into "Novemeber 20" + ", 2023"
$datePara = $tr.Paragraphs(2, 1)
$dateRun = $datePara.Runs(1, 1)

# Insert the new trailing run (keeps the same run formatting) first,
# then shrink the original run's text down to "Novemeber 20".
[void]$dateRun.InsertAfter(", 2023")
$dateRun.Text = "Novemeber 20"
